# Update the header date and every arithmetic-answer cell in the table
# to the values from the target revision. Each old value in this
# document is unique, so a whole-document Find/Replace (MatchWholeWord,
# case-sensitive) on $d.Content reliably retargets exactly one run per
# call, in document order, even where a later replacement's new text
# happens to equal an earlier cell's original text (e.g. "54-31=23").
$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-03-14 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-03-15 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("54-31=23", $true, $false, $false, $false, $false, $true, 1, $false, "31+6=37", 2) | Out-Null
$d.Content.Find.Execute("99+0=99", $true, $false, $false, $false, $false, $true, 1, $false, "72-55=17", 2) | Out-Null
$d.Content.Find.Execute("33-29=4", $true, $false, $false, $false, $false, $true, 1, $false, "76-24=52", 2) | Out-Null
$d.Content.Find.Execute("23+20=43", $true, $false, $false, $false, $false, $true, 1, $false, "42-30=12", 2) | Out-Null
$d.Content.Find.Execute("44+20=64", $true, $false, $false, $false, $false, $true, 1, $false, "66-15=51", 2) | Out-Null
$d.Content.Find.Execute("51+43=94", $true, $false, $false, $false, $false, $true, 1, $false, "89-4=85", 2) | Out-Null
$d.Content.Find.Execute("16-16=0", $true, $false, $false, $false, $false, $true, 1, $false, "57-42=15", 2) | Out-Null
$d.Content.Find.Execute("63-40=23", $true, $false, $false, $false, $false, $true, 1, $false, "57+10=67", 2) | Out-Null
$d.Content.Find.Execute("98-41=57", $true, $false, $false, $false, $false, $true, 1, $false, "29-6=23", 2) | Out-Null
$d.Content.Find.Execute("45-16=29", $true, $false, $false, $false, $false, $true, 1, $false, "78-66=12", 2) | Out-Null
$d.Content.Find.Execute("33-27=6", $true, $false, $false, $false, $false, $true, 1, $false, "5+12=17", 2) | Out-Null
$d.Content.Find.Execute("77-34=43", $true, $false, $false, $false, $false, $true, 1, $false, "23+27=50", 2) | Out-Null
$d.Content.Find.Execute("30-15=15", $true, $false, $false, $false, $false, $true, 1, $false, "40+8=48", 2) | Out-Null
$d.Content.Find.Execute("15+59=74", $true, $false, $false, $false, $false, $true, 1, $false, "77-43=34", 2) | Out-Null
$d.Content.Find.Execute("4+84=88", $true, $false, $false, $false, $false, $true, 1, $false, "57-48=9", 2) | Out-Null
$d.Content.Find.Execute("67+19=86", $true, $false, $false, $false, $false, $true, 1, $false, "13+50=63", 2) | Out-Null
$d.Content.Find.Execute("52+27=79", $true, $false, $false, $false, $false, $true, 1, $false, "68-67=1", 2) | Out-Null
$d.Content.Find.Execute("46+36=82", $true, $false, $false, $false, $false, $true, 1, $false, "4+2=6", 2) | Out-Null
$d.Content.Find.Execute("36+37=73", $true, $false, $false, $false, $false, $true, 1, $false, "90-81=9", 2) | Out-Null
$d.Content.Find.Execute("34+11=45", $true, $false, $false, $false, $false, $true, 1, $false, "31+48=79", 2) | Out-Null
$d.Content.Find.Execute("7+24=31", $true, $false, $false, $false, $false, $true, 1, $false, "50-25=25", 2) | Out-Null
$d.Content.Find.Execute("9+90=99", $true, $false, $false, $false, $false, $true, 1, $false, "40+7=47", 2) | Out-Null
$d.Content.Find.Execute("23+38=61", $true, $false, $false, $false, $false, $true, 1, $false, "45-2=43", 2) | Out-Null
$d.Content.Find.Execute("31+53=84", $true, $false, $false, $false, $false, $true, 1, $false, "4+60=64", 2) | Out-Null
$d.Content.Find.Execute("35+56=91", $true, $false, $false, $false, $false, $true, 1, $false, "34+43=77", 2) | Out-Null
$d.Content.Find.Execute("60+21=81", $true, $false, $false, $false, $false, $true, 1, $false, "22+14=36", 2) | Out-Null
$d.Content.Find.Execute("82+12=94", $true, $false, $false, $false, $false, $true, 1, $false, "12+34=46", 2) | Out-Null
$d.Content.Find.Execute("68+21=89", $true, $false, $false, $false, $false, $true, 1, $false, "39+28=67", 2) | Out-Null
$d.Content.Find.Execute("84-41=43", $true, $false, $false, $false, $false, $true, 1, $false, "2-1=1", 2) | Out-Null
$d.Content.Find.Execute("21+63=84", $true, $false, $false, $false, $false, $true, 1, $false, "96-15=81", 2) | Out-Null
$d.Content.Find.Execute("73+21=94", $true, $false, $false, $false, $false, $true, 1, $false, "77-9=68", 2) | Out-Null
$d.Content.Find.Execute("76+8=84", $true, $false, $false, $false, $false, $true, 1, $false, "66-24=42", 2) | Out-Null
$d.Content.Find.Execute("36+39=75", $true, $false, $false, $false, $false, $true, 1, $false, "91-18=73", 2) | Out-Null
$d.Content.Find.Execute("92-64=28", $true, $false, $false, $false, $false, $true, 1, $false, "53-38=15", 2) | Out-Null
$d.Content.Find.Execute("91+6=97", $true, $false, $false, $false, $false, $true, 1, $false, "0+80=80", 2) | Out-Null
$d.Content.Find.Execute("66-55=11", $true, $false, $false, $false, $false, $true, 1, $false, "40+9=49", 2) | Out-Null
$d.Content.Find.Execute("8+72=80", $true, $false, $false, $false, $false, $true, 1, $false, "68-58=10", 2) | Out-Null
$d.Content.Find.Execute("97-94=3", $true, $false, $false, $false, $false, $true, 1, $false, "42-21=21", 2) | Out-Null
$d.Content.Find.Execute("70-40=30", $true, $false, $false, $false, $false, $true, 1, $false, "46-28=18", 2) | Out-Null
$d.Content.Find.Execute("60-25=35", $true, $false, $false, $false, $false, $true, 1, $false, "69-49=20", 2) | Out-Null
$d.Content.Find.Execute("17+5=22", $true, $false, $false, $false, $false, $true, 1, $false, "35+28=63", 2) | Out-Null
$d.Content.Find.Execute("42-28=14", $true, $false, $false, $false, $false, $true, 1, $false, "28+60=88", 2) | Out-Null
$d.Content.Find.Execute("18+7=25", $true, $false, $false, $false, $false, $true, 1, $false, "32-15=17", 2) | Out-Null
$d.Content.Find.Execute("79-52=27", $true, $false, $false, $false, $false, $true, 1, $false, "87-55=32", 2) | Out-Null
$d.Content.Find.Execute("52-26=26", $true, $false, $false, $false, $false, $true, 1, $false, "2+13=15", 2) | Out-Null
$d.Content.Find.Execute("87-1=86", $true, $false, $false, $false, $false, $true, 1, $false, "98-1=97", 2) | Out-Null
$d.Content.Find.Execute("59-25=34", $true, $false, $false, $false, $false, $true, 1, $false, "21-14=7", 2) | Out-Null
$d.Content.Find.Execute("64+11=75", $true, $false, $false, $false, $false, $true, 1, $false, "19-5=14", 2) | Out-Null
$d.Content.Find.Execute("52+40=92", $true, $false, $false, $false, $false, $true, 1, $false, "73-44=29", 2) | Out-Null
$d.Content.Find.Execute("58-51=7", $true, $false, $false, $false, $false, $true, 1, $false, "54-31=23", 2) | Out-Null
$d.Content.Find.Execute("84+9=93", $true, $false, $false, $false, $false, $true, 1, $false, "63-23=40", 2) | Out-Null
$d.Content.Find.Execute("18-5=13", $true, $false, $false, $false, $false, $true, 1, $false, "99-92=7", 2) | Out-Null
$d.Content.Find.Execute("9+69=78", $true, $false, $false, $false, $false, $true, 1, $false, "57+12=69", 2) | Out-Null
$d.Content.Find.Execute("61+9=70", $true, $false, $false, $false, $false, $true, 1, $false, "77+16=93", 2) | Out-Null
$d.Content.Find.Execute("21+27=48", $true, $false, $false, $false, $false, $true, 1, $false, "9-4=5", 2) | Out-Null
$d.Content.Find.Execute("6+11=17", $true, $false, $false, $false, $false, $true, 1, $false, "8+13=21", 2) | Out-Null
$d.Content.Find.Execute("87-77=10", $true, $false, $false, $false, $false, $true, 1, $false, "82-9=73", 2) | Out-Null
$d.Content.Find.Execute("75+11=86", $true, $false, $false, $false, $false, $true, 1, $false, "46+27=73", 2) | Out-Null
$d.Content.Find.Execute("10+86=96", $true, $false, $false, $false, $false, $true, 1, $false, "28+62=90", 2) | Out-Null
$d.Content.Find.Execute("50-44=6", $true, $false, $false, $false, $false, $true, 1, $false, "32-31=1", 2) | Out-Null
$d.Content.Find.Execute("77-53=24", $true, $false, $false, $false, $false, $true, 1, $false, "79-61=18", 2) | Out-Null
$d.Content.Find.Execute("80-61=19", $true, $false, $false, $false, $false, $true, 1, $false, "70+8=78", 2) | Out-Null
$d.Content.Find.Execute("48-13=35", $true, $false, $false, $false, $false, $true, 1, $false, "0+1=1", 2) | Out-Null
$d.Content.Find.Execute("97-72=25", $true, $false, $false, $false, $false, $true, 1, $false, "92+3=95", 2) | Out-Null
$d.Content.Find.Execute("78+5=83", $true, $false, $false, $false, $false, $true, 1, $false, "79-11=68", 2) | Out-Null
$d.Content.Find.Execute("46-45=1", $true, $false, $false, $false, $false, $true, 1, $false, "92-88=4", 2) | Out-Null
$d.Content.Find.Execute("38+55=93", $true, $false, $false, $false, $false, $true, 1, $false, "9+29=38", 2) | Out-Null
$d.Content.Find.Execute("56-35=21", $true, $false, $false, $false, $false, $true, 1, $false, "58+3=61", 2) | Out-Null
$d.Content.Find.Execute("44+25=69", $true, $false, $false, $false, $false, $true, 1, $false, "29-6=23", 2) | Out-Null
$d.Content.Find.Execute("16+77=93", $true, $false, $false, $false, $false, $true, 1, $false, "34-32=2", 2) | Out-Null
$d.Content.Find.Execute("6+72=78", $true, $false, $false, $false, $false, $true, 1, $false, "43-24=19", 2) | Out-Null
$d.Content.Find.Execute("9+48=57", $true, $false, $false, $false, $false, $true, 1, $false, "53-35=18", 2) | Out-Null
$d.Content.Find.Execute("35+5=40", $true, $false, $false, $false, $false, $true, 1, $false, "75-4=71", 2) | Out-Null
$d.Content.Find.Execute("10+87=97", $true, $false, $false, $false, $false, $true, 1, $false, "1+3=4", 2) | Out-Null
$d.Content.Find.Execute("1+93=94", $true, $false, $false, $false, $false, $true, 1, $false, "19+46=65", 2) | Out-Null
$d.Content.Find.Execute("57-54=3", $true, $false, $false, $false, $false, $true, 1, $false, "18-14=4", 2) | Out-Null
$d.Content.Find.Execute("0+41=41", $true, $false, $false, $false, $false, $true, 1, $false, "79-35=44", 2) | Out-Null
$d.Content.Find.Execute("67+3=70", $true, $false, $false, $false, $false, $true, 1, $false, "11+37=48", 2) | Out-Null
$d.Content.Find.Execute("4+16=20", $true, $false, $false, $false, $false, $true, 1, $false, "21-20=1", 2) | Out-Null
$d.Content.Find.Execute("84-50=34", $true, $false, $false, $false, $false, $true, 1, $false, "29+61=90", 2) | Out-Null
$d.Content.Find.Execute("81-37=44", $true, $false, $false, $false, $false, $true, 1, $false, "22+74=96", 2) | Out-Null
$d.Content.Find.Execute("9+66=75", $true, $false, $false, $false, $false, $true, 1, $false, "37-16=21", 2) | Out-Null
$d.Content.Find.Execute("26-2=24", $true, $false, $false, $false, $false, $true, 1, $false, "35+17=52", 2) | Out-Null
$d.Content.Find.Execute("49-25=24", $true, $false, $false, $false, $false, $true, 1, $false, "17-7=10", 2) | Out-Null
$d.Content.Find.Execute("59-48=11", $true, $false, $false, $false, $false, $true, 1, $false, "2+44=46", 2) | Out-Null
$d.Content.Find.Execute("24+71=95", $true, $false, $false, $false, $false, $true, 1, $false, "84+12=96", 2) | Out-Null
$d.Content.Find.Execute("31-1=30", $true, $false, $false, $false, $false, $true, 1, $false, "79-2=77", 2) | Out-Null
$d.Content.Find.Execute("99-31=68", $true, $false, $false, $false, $false, $true, 1, $false, "2+44=46", 2) | Out-Null
$d.Content.Find.Execute("31-31=0", $true, $false, $false, $false, $false, $true, 1, $false, "32-14=18", 2) | Out-Null
$d.Content.Find.Execute("36+12=48", $true, $false, $false, $false, $false, $true, 1, $false, "59-12=47", 2) | Out-Null
$d.Content.Find.Execute("35+22=57", $true, $false, $false, $false, $false, $true, 1, $false, "74-39=35", 2) | Out-Null
$d.Content.Find.Execute("20-4=16", $true, $false, $false, $false, $false, $true, 1, $false, "61+24=85", 2) | Out-Null
$d.Content.Find.Execute("30+63=93", $true, $false, $false, $false, $false, $true, 1, $false, "96-69=27", 2) | Out-Null
$d.Content.Find.Execute("2+94=96", $true, $false, $false, $false, $false, $true, 1, $false, "18-0=18", 2) | Out-Null
$d.Content.Find.Execute("43+44=87", $true, $false, $false, $false, $false, $true, 1, $false, "27-18=9", 2) | Out-Null
$d.Content.Find.Execute("16+58=74", $true, $false, $false, $false, $false, $true, 1, $false, "61-31=30", 2) | Out-Null
$d.Content.Find.Execute("4+63=67", $true, $false, $false, $false, $false, $true, 1, $false, "73-72=1", 2) | Out-Null
$d.Content.Find.Execute("19+73=92", $true, $false, $false, $false, $false, $true, 1, $false, "1+1=2", 2) | Out-Null
$d.Content.Find.Execute("25-13=12", $true, $false, $false, $false, $false, $true, 1, $false, "56+16=72", 2) | Out-Null
$d.Content.Find.Execute("55-31=24", $true, $false, $false, $false, $false, $true, 1, $false, "80-47=33", 2) | Out-Null
